# Apply "results updated to include safety" changes.
# The set of "reason" categories for each ranked row is recomputed (new
# source data shifted tie-break ordering among equally-ranked reasons),
# while the underlying proportions/ranks (columns D/E) are unaffected.

$wb = $excel.ActiveWorkbook

$wsMoreSafe = $wb.Worksheets.Item("More safe reasons")
$wsLessSafe = $wb.Worksheets.Item("Less safe reasons")

# --- Sheet "More safe reasons" (column C, rows 2-40) ---
$moreSafeValues = @(
    "Increased presence of security personnel",
    "Increased awareness about health risks and appropriate behaviour",
    "Introduction of an evening curfew",
    "Implementation of rules and regulations",
    "Markets less crowded",
    "Relaxation of rules and regulations",
    "Introduction of new rules and regulations",
    "No theft",
    "n/a",
    "Installation of sanitary infrastructure",
    "Appropriate behaviour of security personnel",
    "Reduced threats of disease spread",
    "No threats of covic infections",
    "Increased presence of security personnel",
    "Implementation of rules and regulations",
    "Introduction of an evening curfew",
    "Increased awareness about health risks and appropriate behaviour",
    "Introduction of new rules and regulations",
    "No threats of covic infections",
    "Relaxation of rules and regulations",
    "Reduced threats of disease spread",
    "Installation of sanitary infrastructure",
    "n/a",
    "No theft",
    "Appropriate behaviour of security personnel",
    "Markets less crowded",
    "Increased presence of security personnel",
    "Increased awareness about health risks and appropriate behaviour",
    "Introduction of an evening curfew",
    "Implementation of rules and regulations",
    "Markets less crowded",
    "No theft",
    "Relaxation of rules and regulations",
    "n/a",
    "Installation of sanitary infrastructure",
    "Introduction of new rules and regulations",
    "Appropriate behaviour of security personnel",
    "No threats of covic infections",
    "Reduced threats of disease spread"
)
for ($i = 0; $i -lt $moreSafeValues.Count; $i++) {
    $row = $i + 2
    $wsMoreSafe.Range("C$row").Value = $moreSafeValues[$i]
}

# --- Sheet "Less safe reasons" (column C, rows 2-34) ---
$lessSafeValues = @(
    "Fear of theft ",
    "Non-compliance with official rules and regulations",
    "Fear of contracting COVID-19",
    "Discretionary behaviour of security personnel",
    "Health risks in relation to COVID-19",
    "Corona-infected individuals in the area",
    "Fear of contracting COVID-20",
    "Some people stopped working and may pose a big threat to our merchandise",
    "Lack of customers",
    "Restrictive rules and regulations",
    "Increased presence of security personnel",
    "Fear of theft ",
    "Fear of contracting COVID-19",
    "Non-compliance with official rules and regulations",
    "Some people stopped working and may pose a big threat to our merchandise",
    "Corona-infected individuals in the area",
    "Fear of contracting COVID-20",
    "Health risks in relation to COVID-19",
    "Restrictive rules and regulations",
    "Lack of customers",
    "Discretionary behaviour of security personnel",
    "Increased presence of security personnel",
    "Non-compliance with official rules and regulations",
    "Fear of theft ",
    "Fear of contracting COVID-19",
    "Health risks in relation to COVID-19",
    "Discretionary behaviour of security personnel",
    "Lack of customers",
    "Fear of contracting COVID-20",
    "Increased presence of security personnel",
    "Corona-infected individuals in the area",
    "Some people stopped working and may pose a big threat to our merchandise",
    "Restrictive rules and regulations"
)
for ($i = 0; $i -lt $lessSafeValues.Count; $i++) {
    $row = $i + 2
    $wsLessSafe.Range("C$row").Value = $lessSafeValues[$i]
}
